$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Address, $Text)
    $ws.Range($Address).NumberFormat = "@"
    $ws.Range($Address).Value = $Text
    $ws.Range($Address).Style = "Normal"
}

$ws.Range("D2").Value = "42.725.88"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.547.69"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextCell "D5" "310.03"
$ws.Range("E5").Value = "  -2.50%  "
Set-TextCell "D6" "99.38"
$ws.Range("E6").Value = "  +2.38%  "
Set-TextCell "D7" "0.571"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextCell "D9" "0.533"
$ws.Range("E9").Value = "  -0.38%  "
Set-TextCell "D10" "35.86"
$ws.Range("E10").Value = "  -0.72%  "
Set-TextCell "D11" "0.0807"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "2.934.85"
$ws.Range("E14").Value = "  +0.34%  "
Set-TextCell "D15" "16.01"
$ws.Range("E15").Value = "  +5.53%  "
$ws.Range("D16").Value = "2.562.19"
$ws.Range("E16").Value = "  -1.03%  "
Set-TextCell "D17" "0.841"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "42.729.07"
$ws.Range("E18").Value = "  -0.51%  "
Set-TextCell "D19" "6.76"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D20" "12.39"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0956"
$ws.Range("E21").Value = "  -1.26%  "
Set-TextCell "D22" "69.48"
$ws.Range("E22").Value = "  -0.38%  "
Set-TextCell "D23" "248.56"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("E25").Value = "  -0.36%  "
Set-TextCell "D26" "26.60"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -1.32%  "
Set-TextCell "D29" "40.26"
$ws.Range("E29").Value = "  -1.53%  "
Set-TextCell "D30" "10.12"
$ws.Range("E30").Value = "  -3.80%  "
Set-TextCell "D31" "157.92"
$ws.Range("E31").Value = "  +0.18%  "
Set-TextCell "D32" "5.74"
$ws.Range("E32").Value = "  -2.83%  "
Set-TextCell "D33" "0.0807"
$ws.Range("E33").Value = "  +1.61%  "
Set-TextCell "D34" "3.30"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  -3.37%  "
Set-TextCell "D36" "2.63"
$ws.Range("E36").Value = "  -3.31%  "
Set-TextCell "D37" "2.61"
$ws.Range("E37").Value = "  +5.74%  "
Set-TextCell "D38" "18.24"
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("E39").Value = "  -1.37%  "
Set-TextCell "D40" "0.118"
$ws.Range("E40").Value = "  -0.77%  "
Set-TextCell "D41" "22.55"
$ws.Range("E41").Value = "  +1.80%  "
Set-TextCell "D42" "4.19"
$ws.Range("E42").Value = "  +9.53%  "
$ws.Range("E43").Value = "  -0.14%  "
Set-TextCell "D44" "0.0301"
$ws.Range("E44").Value = "  -1.47%  "
Set-TextCell "D45" "3.26"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "1.988.02"
$ws.Range("E46").Value = "  -1.15%  "
Set-TextCell "D47" "9.02"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").Value = "2.790.55"
$ws.Range("E48").Value = "  +0.34%  "
Set-TextCell "D49" "81.33"
$ws.Range("E49").Value = "  -3.48%  "
Set-TextCell "D50" "0.194"
$ws.Range("E50").Value = "  +1.04%  "
Set-TextCell "D51" "73.40"
$ws.Range("E51").Value = "  -3.06%  "
